$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("procedimientos")

# Fix the casing of the procedure names in column B (rows 2 and 3)
$ws.Range("B2").Value = "val_interfaz_a01"
$ws.Range("B3").Value = "val_interfaz_b01"
